$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 6873.269  # H69: 6591.7666 -> 6873.269
$ws.Cells.Item(69, 9).Value = 6450  # I69: 0 -> 6450
$ws.Cells.Item(69, 10).Value = 6890.2  # J69: 6591.7666 -> 6890.2
$ws.Cells.Item(69, 11).Value = 19350  # K69: 0 -> 19350
$ws.Cells.Item(69, 12).Value = 20670.6  # L69: 19775.2998 -> 20670.6
$ws.Cells.Item(69, 13).Value = -18476  # M69: None -> -18476
$ws.Cells.Item(69, 14).Value = -22418.6  # N69: -21523.2998 -> -22418.6
$ws.Cells.Item(72, 8).Value = 6873.269  # H72: 6591.7666 -> 6873.269
$ws.Cells.Item(72, 9).Value = 6450  # I72: 0 -> 6450
$ws.Cells.Item(72, 10).Value = 6890.2  # J72: 6591.7666 -> 6890.2
$ws.Cells.Item(72, 11).Value = 58050  # K72: 0 -> 58050
$ws.Cells.Item(72, 12).Value = 62011.8  # L72: 59325.8994 -> 62011.8
$ws.Cells.Item(72, 13).Value = -53682  # M72: None -> -53682
$ws.Cells.Item(72, 14).Value = -70747.79999999999  # N72: -68061.89939999999 -> -70747.79999999999
$ws.Cells.Item(74, 8).Value = 5766.6904  # H74: 5482.2446 -> 5766.6904
$ws.Cells.Item(74, 9).Value = 3811.625  # I74: 3446.6316 -> 3811.625
$ws.Cells.Item(74, 11).Value = 3811.625  # K74: 3446.6316 -> 3811.625
$ws.Cells.Item(74, 13).Value = -2875.625  # M74: -2510.6316 -> -2875.625
$ws.Cells.Item(77, 8).Value = 5766.6904  # H77: 5482.2446 -> 5766.6904
$ws.Cells.Item(77, 9).Value = 3811.625  # I77: 3446.6316 -> 3811.625
$ws.Cells.Item(77, 11).Value = 19058.125  # K77: 17233.158 -> 19058.125
$ws.Cells.Item(77, 13).Value = -14378.125  # M77: -12553.158 -> -14378.125
$ws.Cells.Item(112, 8).Value = 6229.857  # H112: 6001.273 -> 6229.857
$ws.Cells.Item(112, 10).Value = 7889.9375  # J112: 7496.4707 -> 7889.9375
$ws.Cells.Item(112, 12).Value = 23669.8125  # L112: 22489.4121 -> 23669.8125
$ws.Cells.Item(112, 14).Value = -25885.8125  # N112: -24705.4121 -> -25885.8125
$ws.Cells.Item(133, 8).Value = 0  # H133: 86666.664 -> 0
$ws.Cells.Item(133, 10).Value = 0  # J133: 86666.664 -> 0
$ws.Cells.Item(133, 12).Value = 0  # L133: 86666.664 -> 0
$ws.Cells.Item(133, 14).ClearContents()  # N133: -96786.664 -> (blank)
$ws.Cells.Item(137, 8).Value = 33560.395  # H137: 33622.02 -> 33560.395
$ws.Cells.Item(137, 9).Value = 45745.5  # I137: 46937.9 -> 45745.5
$ws.Cells.Item(137, 10).Value = 3097.625  # J137: 3073.8235 -> 3097.625
$ws.Cells.Item(137, 11).Value = 137236.5  # K137: 140813.7 -> 137236.5
$ws.Cells.Item(137, 12).Value = 9292.875  # L137: 9221.470499999999 -> 9292.875
$ws.Cells.Item(137, 13).Value = -134686.5  # M137: -138263.7 -> -134686.5
$ws.Cells.Item(137, 14).Value = -14392.875  # N137: -14321.4705 -> -14392.875
$ws.Cells.Item(138, 8).Value = 3152.7104  # H138: 3081 -> 3152.7104
$ws.Cells.Item(138, 10).Value = 3893.318  # J138: 3663.5386 -> 3893.318
$ws.Cells.Item(138, 12).Value = 11679.954  # L138: 10990.6158 -> 11679.954
$ws.Cells.Item(138, 14).Value = -21959.954  # N138: -21270.6158 -> -21959.954

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7487.53  # H32: 8509.736999999999 -> 7487.53
$ws.Cells.Item(32, 9).Value = 3822.6233  # I32: 5063.5127 -> 3822.6233
$ws.Cells.Item(32, 10).Value = 19757  # J32: 21310 -> 19757
$ws.Cells.Item(32, 11).Value = 3822.6233  # K32: 5063.5127 -> 3822.6233
$ws.Cells.Item(32, 12).Value = 19757  # L32: 21310 -> 19757
$ws.Cells.Item(32, 13).Value = -3535.6233  # M32: -4776.5127 -> -3535.6233
$ws.Cells.Item(32, 14).Value = -20331  # N32: -21884 -> -20331
$ws.Cells.Item(45, 8).Value = 7146264.5  # H45: 6806035.5 -> 7146264.5
$ws.Cells.Item(45, 9).Value = 10990639  # I45: 10205697 -> 10990639
$ws.Cells.Item(45, 11).Value = 10990639  # K45: 10205697 -> 10990639
$ws.Cells.Item(45, 13).Value = -10990262  # M45: -10205320 -> -10990262
$ws.Cells.Item(74, 8).Value = 28461.97  # H74: 29268.559 -> 28461.97
$ws.Cells.Item(74, 10).Value = 103946.664  # J74: 116810.25 -> 103946.664
$ws.Cells.Item(74, 12).Value = 103946.664  # L74: 116810.25 -> 103946.664
$ws.Cells.Item(74, 14).Value = -105694.664  # N74: -118558.25 -> -105694.664
$ws.Cells.Item(77, 8).Value = 28461.97  # H77: 29268.559 -> 28461.97
$ws.Cells.Item(77, 10).Value = 103946.664  # J77: 116810.25 -> 103946.664
$ws.Cells.Item(77, 12).Value = 519733.32  # L77: 584051.25 -> 519733.32
$ws.Cells.Item(77, 14).Value = -528469.3200000001  # N77: -592787.25 -> -528469.3200000001
$ws.Cells.Item(80, 8).Value = 46499.5  # H80: 60000 -> 46499.5
$ws.Cells.Item(80, 9).Value = 46499.5  # I80: 60000 -> 46499.5
$ws.Cells.Item(80, 11).Value = 46499.5  # K80: 60000 -> 46499.5
$ws.Cells.Item(80, 13).Value = -45501.5  # M80: -59002 -> -45501.5
$ws.Cells.Item(83, 8).Value = 46499.5  # H83: 60000 -> 46499.5
$ws.Cells.Item(83, 9).Value = 46499.5  # I83: 60000 -> 46499.5
$ws.Cells.Item(83, 11).Value = 139498.5  # K83: 180000 -> 139498.5
$ws.Cells.Item(83, 13).Value = -134506.5  # M83: -175008 -> -134506.5
$ws.Cells.Item(102, 8).Value = 2978299.5  # H102: 2978300.2 -> 2978299.5
$ws.Cells.Item(102, 9).Value = 3088421.8  # I102: 3088422.5 -> 3088421.8
$ws.Cells.Item(102, 11).Value = 3088421.8  # K102: 3088422.5 -> 3088421.8
$ws.Cells.Item(102, 13).Value = -3086799.8  # M102: -3086800.5 -> -3086799.8
$ws.Cells.Item(110, 8).Value = 1737655.4  # H110: 2138391.2 -> 1737655.4
$ws.Cells.Item(110, 9).Value = 2138191.5  # I110: 2526826.2 -> 2138191.5
$ws.Cells.Item(110, 11).Value = 2138191.5  # K110: 2526826.2 -> 2138191.5
$ws.Cells.Item(110, 13).Value = -2136146.5  # M110: -2524781.2 -> -2136146.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 3666.6667  # H7: 3000 -> 3666.6667
$ws.Cells.Item(86, 8).Value = 25493704  # H86: 27086974 -> 25493704
$ws.Cells.Item(86, 9).Value = 48150924  # I86: 54169616 -> 48150924
$ws.Cells.Item(86, 11).Value = 48150924  # K86: 54169616 -> 48150924
$ws.Cells.Item(86, 13).Value = -48149801  # M86: -54168493 -> -48149801
$ws.Cells.Item(89, 8).Value = 25493704  # H89: 27086974 -> 25493704
$ws.Cells.Item(89, 9).Value = 48150924  # I89: 54169616 -> 48150924
$ws.Cells.Item(89, 11).Value = 240754620  # K89: 270848080 -> 240754620
$ws.Cells.Item(89, 13).Value = -240749004  # M89: -270842464 -> -240749004
$ws.Cells.Item(94, 8).Value = 4356100  # H94: 4554091 -> 4356100
$ws.Cells.Item(94, 9).Value = 4763395.5  # I94: 5001550 -> 4763395.5
$ws.Cells.Item(94, 11).Value = 4763395.5  # K94: 5001550 -> 4763395.5
$ws.Cells.Item(94, 13).Value = -4762944.5  # M94: -5001099 -> -4762944.5
$ws.Cells.Item(105, 8).Value = 31251000  # H105: 12507001 -> 31251000
$ws.Cells.Item(105, 9).Value = 31251000  # I105: 15631251 -> 31251000
$ws.Cells.Item(105, 10).Value = 0  # J105: 9999 -> 0
$ws.Cells.Item(105, 11).Value = 31251000  # K105: 15631251 -> 31251000
$ws.Cells.Item(105, 12).Value = 0  # L105: 9999 -> 0
$ws.Cells.Item(105, 13).Value = -31249253  # M105: -15629504 -> -31249253
$ws.Cells.Item(105, 14).ClearContents()  # N105: -13493 -> (blank)
$ws.Cells.Item(107, 8).Value = 5955515.5  # H107: 7940537 -> 5955515.5
$ws.Cells.Item(107, 9).Value = 8929543  # I107: 14286998 -> 8929543
$ws.Cells.Item(107, 11).Value = 8929543  # K107: 14286998 -> 8929543
$ws.Cells.Item(107, 13).Value = -8927623  # M107: -14285078 -> -8927623
$ws.Cells.Item(134, 8).Value = 3592.8262  # H134: 3668.6 -> 3592.8262
$ws.Cells.Item(134, 9).Value = 1044.3158  # I134: 1067.5946 -> 1044.3158
$ws.Cells.Item(134, 11).Value = 3132.9474  # K134: 3202.7838 -> 3132.9474
$ws.Cells.Item(134, 13).Value = -597.9474  # M134: -667.7837999999997 -> -597.9474

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 29560.19  # H31: 15118.571 -> 29560.19
$ws.Cells.Item(31, 9).Value = 0  # I31: 1847.8286 -> 0
$ws.Cells.Item(31, 10).Value = 29560.19  # J31: 26177.523 -> 29560.19
$ws.Cells.Item(31, 11).Value = 0  # K31: 1847.8286 -> 0
$ws.Cells.Item(31, 12).Value = 29560.19  # L31: 26177.523 -> 29560.19
$ws.Cells.Item(31, 13).ClearContents()  # M31: -1552.8286 -> (blank)
$ws.Cells.Item(31, 14).Value = -30150.19  # N31: -26767.523 -> -30150.19
$ws.Cells.Item(34, 8).Value = 29560.19  # H34: 15118.571 -> 29560.19
$ws.Cells.Item(34, 9).Value = 0  # I34: 1847.8286 -> 0
$ws.Cells.Item(34, 10).Value = 29560.19  # J34: 26177.523 -> 29560.19
$ws.Cells.Item(34, 11).Value = 0  # K34: 1847.8286 -> 0
$ws.Cells.Item(34, 12).Value = 29560.19  # L34: 26177.523 -> 29560.19
$ws.Cells.Item(34, 13).ClearContents()  # M34: -1645.8286 -> (blank)
$ws.Cells.Item(34, 14).Value = -29964.19  # N34: -26581.523 -> -29964.19
$ws.Cells.Item(58, 8).Value = 8168.9614  # H58: 5295.375 -> 8168.9614
$ws.Cells.Item(58, 9).Value = 12922.462  # I58: 6581.6313 -> 12922.462
$ws.Cells.Item(58, 11).Value = 12922.462  # K58: 6581.6313 -> 12922.462
$ws.Cells.Item(58, 13).Value = -12719.462  # M58: -6378.6313 -> -12719.462
$ws.Cells.Item(62, 8).Value = 3908.5715  # H62: 4589 -> 3908.5715
$ws.Cells.Item(62, 9).Value = 4587.5  # I62: 5216.6665 -> 4587.5
$ws.Cells.Item(62, 10).Value = 3003.3333  # J62: 3647.5 -> 3003.3333
$ws.Cells.Item(62, 11).Value = 4587.5  # K62: 5216.6665 -> 4587.5
$ws.Cells.Item(62, 12).Value = 3003.3333  # L62: 3647.5 -> 3003.3333
$ws.Cells.Item(62, 13).Value = -3963.5  # M62: -4592.6665 -> -3963.5
$ws.Cells.Item(62, 14).Value = -4251.3333  # N62: -4895.5 -> -4251.3333
$ws.Cells.Item(65, 8).Value = 3908.5715  # H65: 4589 -> 3908.5715
$ws.Cells.Item(65, 9).Value = 4587.5  # I65: 5216.6665 -> 4587.5
$ws.Cells.Item(65, 10).Value = 3003.3333  # J65: 3647.5 -> 3003.3333
$ws.Cells.Item(65, 11).Value = 22937.5  # K65: 26083.3325 -> 22937.5
$ws.Cells.Item(65, 12).Value = 15016.6665  # L65: 18237.5 -> 15016.6665
$ws.Cells.Item(65, 13).Value = -19817.5  # M65: -22963.3325 -> -19817.5
$ws.Cells.Item(65, 14).Value = -21256.6665  # N65: -24477.5 -> -21256.6665
$ws.Cells.Item(132, 8).Value = 47129.668  # H132: 44042.89 -> 47129.668
$ws.Cells.Item(132, 9).Value = 31943.516  # I132: 29356.25 -> 31943.516
$ws.Cells.Item(132, 10).Value = 102812.22  # J132: 102789.445 -> 102812.22
$ws.Cells.Item(132, 11).Value = 95830.548  # K132: 88068.75 -> 95830.548
$ws.Cells.Item(132, 12).Value = 308436.66  # L132: 308368.335 -> 308436.66
$ws.Cells.Item(132, 13).Value = -93300.548  # M132: -85538.75 -> -93300.548
$ws.Cells.Item(132, 14).Value = -313496.66  # N132: -313428.335 -> -313496.66
$ws.Cells.Item(134, 8).Value = 2255.8538  # H134: 2308.875 -> 2255.8538
$ws.Cells.Item(134, 9).Value = 1311.5807  # I134: 1350.8 -> 1311.5807
$ws.Cells.Item(134, 11).Value = 3934.7421  # K134: 4052.4 -> 3934.7421
$ws.Cells.Item(134, 13).Value = -1399.7421  # M134: -1517.4 -> -1399.7421
$ws.Cells.Item(136, 8).Value = 8168.9614  # H136: 5295.375 -> 8168.9614
$ws.Cells.Item(136, 9).Value = 12922.462  # I136: 6581.6313 -> 12922.462
$ws.Cells.Item(136, 11).Value = 38767.386  # K136: 19744.8939 -> 38767.386
$ws.Cells.Item(136, 13).Value = -36217.386  # M136: -17194.8939 -> -36217.386

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 436177.53  # H2: 467328.75 -> 436177.53
$ws.Cells.Item(2, 9).Value = 575.6667  # I2: 577.3333 -> 575.6667
$ws.Cells.Item(2, 10).Value = 817329.2  # J2: 934080.1 -> 817329.2
$ws.Cells.Item(2, 11).Value = 3454.0002  # K2: 3463.9998 -> 3454.0002
$ws.Cells.Item(2, 12).Value = 4903975.199999999  # L2: 5604480.6 -> 4903975.199999999
$ws.Cells.Item(2, 13).Value = -3341.0002  # M2: -3350.9998 -> -3341.0002
$ws.Cells.Item(2, 14).Value = -4904201.199999999  # N2: -5604706.6 -> -4904201.199999999
$ws.Cells.Item(4, 8).Value = 58577430  # H4: 46861960 -> 58577430
$ws.Cells.Item(4, 9).Value = 72128010  # I4: 57702416 -> 72128010
$ws.Cells.Item(4, 10).Value = 17925700  # J4: 14340581 -> 17925700
$ws.Cells.Item(4, 11).Value = 216384030  # K4: 173107248 -> 216384030
$ws.Cells.Item(4, 12).Value = 53777100  # L4: 43021743 -> 53777100
$ws.Cells.Item(4, 13).Value = -216383918  # M4: -173107136 -> -216383918
$ws.Cells.Item(4, 14).Value = -53777324  # N4: -43021967 -> -53777324
$ws.Cells.Item(133, 8).Value = 2083.6667  # H133: 2129 -> 2083.6667
$ws.Cells.Item(133, 9).Value = 2083.6667  # I133: 2129 -> 2083.6667
$ws.Cells.Item(133, 11).Value = 6251.000100000001  # K133: 6387 -> 6251.000100000001
$ws.Cells.Item(133, 13).Value = -1191.000100000001  # M133: -1327 -> -1191.000100000001
$ws.Cells.Item(137, 8).Value = 2153.5  # H137: 1685.1875 -> 2153.5
$ws.Cells.Item(137, 9).Value = 1803.909  # I137: 1397.6 -> 1803.909
$ws.Cells.Item(137, 11).Value = 5411.727000000001  # K137: 4192.799999999999 -> 5411.727000000001
$ws.Cells.Item(137, 13).Value = -311.7270000000008  # M137: 907.2000000000007 -> -311.7270000000008

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 33354834  # H70: 33355134 -> 33354834
$ws.Cells.Item(70, 9).Value = 50003100  # I70: 50003550 -> 50003100
$ws.Cells.Item(70, 11).Value = 50003100  # K70: 50003550 -> 50003100
$ws.Cells.Item(70, 13).Value = -50002830  # M70: -50003280 -> -50002830
$ws.Cells.Item(73, 8).Value = 33354834  # H73: 33355134 -> 33354834
$ws.Cells.Item(73, 9).Value = 50003100  # I73: 50003550 -> 50003100
$ws.Cells.Item(73, 11).Value = 50003100  # K73: 50003550 -> 50003100
$ws.Cells.Item(73, 13).Value = -50002164  # M73: -50002614 -> -50002164
$ws.Cells.Item(124, 8).Value = 0  # H124: 62000 -> 0
$ws.Cells.Item(124, 10).Value = 0  # J124: 62000 -> 0
$ws.Cells.Item(124, 12).Value = 0  # L124: 62000 -> 0
$ws.Cells.Item(124, 14).ClearContents()  # N124: -71820 -> (blank)
$ws.Cells.Item(132, 8).Value = 2259.6338  # H132: 2307.7537 -> 2259.6338
$ws.Cells.Item(132, 9).Value = 2105.3403  # I132: 2172.2666 -> 2105.3403
$ws.Cells.Item(132, 11).Value = 6316.0209  # K132: 6516.7998 -> 6316.0209
$ws.Cells.Item(132, 13).Value = -3786.0209  # M132: -3986.7998 -> -3786.0209

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 24749.125  # H17: 0 -> 24749.125
$ws.Cells.Item(17, 9).Value = 19999  # I17: 0 -> 19999
$ws.Cells.Item(17, 10).Value = 38999.5  # J17: 0 -> 38999.5
$ws.Cells.Item(17, 11).Value = 19999  # K17: 0 -> 19999
$ws.Cells.Item(17, 12).Value = 38999.5  # L17: 0 -> 38999.5
$ws.Cells.Item(17, 13).Value = -19829  # M17: None -> -19829
$ws.Cells.Item(17, 14).Value = -39339.5  # N17: None -> -39339.5
$ws.Cells.Item(22, 8).Value = 39461.625  # H22: 85262.73 -> 39461.625
$ws.Cells.Item(22, 9).Value = 89728  # I22: 888888 -> 89728
$ws.Cells.Item(22, 10).Value = 3557.0715  # J22: 4900.2 -> 3557.0715
$ws.Cells.Item(22, 11).Value = 89728  # K22: 888888 -> 89728
$ws.Cells.Item(22, 12).Value = 3557.0715  # L22: 4900.2 -> 3557.0715
$ws.Cells.Item(22, 13).Value = -89433  # M22: -888593 -> -89433
$ws.Cells.Item(22, 14).Value = -4147.0715  # N22: -5490.2 -> -4147.0715
$ws.Cells.Item(25, 8).Value = 23249  # H25: 1500 -> 23249
$ws.Cells.Item(25, 9).Value = 19999  # I25: 1500 -> 19999
$ws.Cells.Item(25, 10).Value = 39499  # J25: 0 -> 39499
$ws.Cells.Item(25, 11).Value = 19999  # K25: 1500 -> 19999
$ws.Cells.Item(25, 12).Value = 39499  # L25: 0 -> 39499
$ws.Cells.Item(25, 13).Value = -19769  # M25: -1270 -> -19769
$ws.Cells.Item(25, 14).Value = -39959  # N25: None -> -39959
$ws.Cells.Item(27, 8).Value = 39461.625  # H27: 85262.73 -> 39461.625
$ws.Cells.Item(27, 9).Value = 89728  # I27: 888888 -> 89728
$ws.Cells.Item(27, 10).Value = 3557.0715  # J27: 4900.2 -> 3557.0715
$ws.Cells.Item(27, 11).Value = 89728  # K27: 888888 -> 89728
$ws.Cells.Item(27, 12).Value = 3557.0715  # L27: 4900.2 -> 3557.0715
$ws.Cells.Item(27, 13).Value = -89621  # M27: -888781 -> -89621
$ws.Cells.Item(27, 14).Value = -3771.0715  # N27: -5114.2 -> -3771.0715
$ws.Cells.Item(31, 8).Value = 3858.2222  # H31: 4572.2 -> 3858.2222
$ws.Cells.Item(31, 9).Value = 4683  # I31: 2428.5 -> 4683
$ws.Cells.Item(31, 10).Value = 3198.4  # J31: 6001.3335 -> 3198.4
$ws.Cells.Item(31, 11).Value = 4683  # K31: 2428.5 -> 4683
$ws.Cells.Item(31, 12).Value = 3198.4  # L31: 6001.3335 -> 3198.4
$ws.Cells.Item(31, 13).Value = -4435  # M31: -2180.5 -> -4435
$ws.Cells.Item(31, 14).Value = -3694.4  # N31: -6497.3335 -> -3694.4
$ws.Cells.Item(46, 8).Value = 4522  # H46: 4523.905 -> 4522
$ws.Cells.Item(46, 9).Value = 1517.2858  # I46: 1528.7142 -> 1517.2858
$ws.Cells.Item(46, 11).Value = 1517.2858  # K46: 1528.7142 -> 1517.2858
$ws.Cells.Item(46, 13).Value = -1329.2858  # M46: -1340.7142 -> -1329.2858
$ws.Cells.Item(55, 8).Value = 2137.389  # H55: 2040.6842 -> 2137.389
$ws.Cells.Item(55, 9).Value = 2659.8333  # I55: 2322.7144 -> 2659.8333
$ws.Cells.Item(55, 11).Value = 2659.8333  # K55: 2322.7144 -> 2659.8333
$ws.Cells.Item(55, 13).Value = -2486.8333  # M55: -2149.7144 -> -2486.8333
$ws.Cells.Item(132, 8).Value = 7276  # H132: 7360.775 -> 7276
$ws.Cells.Item(132, 9).Value = 7425.6387  # I132: 7428.6943 -> 7425.6387
$ws.Cells.Item(132, 10).Value = 6198.6  # J132: 6749.5 -> 6198.6
$ws.Cells.Item(132, 11).Value = 22276.9161  # K132: 22286.0829 -> 22276.9161
$ws.Cells.Item(132, 12).Value = 18595.8  # L132: 20248.5 -> 18595.8
$ws.Cells.Item(132, 13).Value = -19746.9161  # M132: -19756.0829 -> -19746.9161
$ws.Cells.Item(132, 14).Value = -23655.8  # N132: -25308.5 -> -23655.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 84475  # H46: 70285.8 -> 84475
$ws.Cells.Item(46, 10).Value = 84475  # J46: 70285.8 -> 84475
$ws.Cells.Item(46, 12).Value = 84475  # L46: 70285.8 -> 84475
$ws.Cells.Item(46, 14).Value = -84937  # N46: -70747.8 -> -84937
$ws.Cells.Item(113, 8).Value = 861.619  # H113: 900.2 -> 861.619
$ws.Cells.Item(113, 10).Value = 1747.125  # J113: 1983.8572 -> 1747.125
$ws.Cells.Item(113, 12).Value = 5241.375  # L113: 5951.571599999999 -> 5241.375
$ws.Cells.Item(113, 14).Value = -9581.375  # N113: -10291.5716 -> -9581.375
$ws.Cells.Item(115, 8).Value = 46000  # H115: 35333.332 -> 46000
$ws.Cells.Item(115, 10).Value = 46000  # J115: 35333.332 -> 46000
$ws.Cells.Item(115, 12).Value = 46000  # L115: 35333.332 -> 46000
$ws.Cells.Item(115, 14).Value = -49134  # N115: -38467.332 -> -49134
$ws.Cells.Item(134, 8).Value = 84475  # H134: 70285.8 -> 84475
$ws.Cells.Item(134, 10).Value = 84475  # J134: 70285.8 -> 84475
$ws.Cells.Item(134, 12).Value = 253425  # L134: 210857.4 -> 253425
$ws.Cells.Item(134, 14).Value = -258495  # N134: -215927.4 -> -258495
